$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jobs")

# Insert a new row for the "beginner" job at row 2
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "beginner"
$ws.Range("C2").Value = "no job"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 10

# Add new "Enabled" column header (copy header style from A1)
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Enabled"

# Set Enabled values for each job row
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $true
$ws.Range("F4").Value = $true
$ws.Range("F5").Value = $true
$ws.Range("F6").Value = $false
$ws.Range("F7").Value = $false
$ws.Range("F8").Value = $false
$ws.Range("F9").Value = $true
$ws.Range("F10").Value = $true
$ws.Range("F11").Value = $true
$ws.Range("F12").Value = $true

# Match the final selection state left behind in the saved file
$ws.Range("F16").Select() | Out-Null
